$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "line" series grows from line1..line6 to line1..line8, which pushes the
# "extr" series (extr1..extr8) down by two rows. Rows 8-9 (previously the
# start of the extr series) now become line7/line8, rows 10-15 become
# extr1..extr6, and two brand new rows (16-17) are appended for extr7/extr8.

$rows = @(
    @{ Row = 8;  Name = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  Name = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; Name = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; Name = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; Name = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; Name = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; Name = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; Name = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; Name = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; Name = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
}

# New rows 16 and 17 need the A column populated (sequence continues 14, 15)
# with the same formatting as the rest of column A (bold, bordered, centered) -
# copy the format straight from the row above instead of rebuilding it.
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
